# Sync attendance_reports: normalize "Recorded By" (column G) entries so that
# an entry ending with the literal token "System" has that token moved to the
# front of the comma-separated list (the displaced token keeps its original
# casing and moves to the end of the list). Entries that already start with
# "System", are a single token, or involve "admin@admin.com" are left as-is.

$excel = New-Object -ComObject Excel.Application
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null) {
        $parts = $val -split ", "
        $n = $parts.Count

        if ($n -gt 1) {
            $firstTok = $parts[0]
            $lastTok = $parts[$n - 1]
            $hasAdmin = $val.Contains("admin@admin.com")

            if ($lastTok.Equals("System") -and (-not $firstTok.Equals("System")) -and (-not $hasAdmin)) {
                $parts[0] = $lastTok
                $parts[$n - 1] = $firstTok
                $newVal = $parts -join ", "
                $cell.Value2 = $newVal
            }
        }
    }
}
